# feat: add 2022-Q1 data
#
# Inserts a new worksheet "2022-Q1" right before the "总计" (total) sheet,
# populates it with the Q1-2022 fund-holding data, and prepends a matching
# summary row ("2022-Q1", 5, 1.73) to the "总计" sheet.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# --- 1) Create the new "2022-Q1" sheet, placed before "总计" -----------------
$new = $wb.Worksheets.Add($total)
$new.Name = "2022-Q1"

# Reuse the "2021-Q4" sheet's look (bold/centered/bordered header row + bold
# column A) by copying its formatting onto the new sheet. (A1 is left alone,
# same as the reference sheet, which has no A1 cell at all.)
$q4.Range("B1:H1").Copy()
$new.Range("B1").PasteSpecial(-4122)
$q4.Range("A2:H2").Copy()
$new.Range("A2:H6").PasteSpecial(-4122)

# --- 2) Header row -----------------------------------------------------------
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# Columns B (fund codes, e.g. "005368" - leading zeros must survive) and
# D:G (numeric-looking figures like "23.52", "88.60") must stay as TEXT
# rather than being auto-coerced into numbers.
$new.Range("B2:B6").NumberFormat = "@"
$new.Range("D2:G6").NumberFormat = "@"

# --- 3) Data rows --------------------------------------------------------------
$data = @(
    @(0, "005368", "富国清洁能源产业灵活配置混合A", "23.52", "88.60", "3.06", "0.7197", 7),
    @(1, "001556", "天弘中证500指数增强A",           "41.41", "94.29", "1.71", "0.7081", 8),
    @(2, "001557", "天弘中证500指数增强C",           "13.97", "94.29", "1.71", "0.2389", 8),
    @(3, "011127", "富国清洁能源产业灵活配置混合C", "1.61",  "88.60", "3.06", "0.0493", 7),
    @(4, "159962", "华夏中证四川国企改革ETF",        "0.49",  "95.82", "3.26", "0.0160", 6)
)

$r = 2
foreach ($row in $data) {
    $new.Cells.Item($r, 1).Value = $row[0]
    $new.Cells.Item($r, 2).Value = $row[1]
    $new.Cells.Item($r, 3).Value = $row[2]
    $new.Cells.Item($r, 4).Value = $row[3]
    $new.Cells.Item($r, 5).Value = $row[4]
    $new.Cells.Item($r, 6).Value = $row[5]
    $new.Cells.Item($r, 7).Value = $row[6]
    $new.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# The "@" text NumberFormat above was only needed transiently so Excel didn't
# coerce the numeric-looking strings into real numbers; drop it now so the
# cells end up as plain (unstyled) text, same as columns C/H.
$new.Range("B2:B6").ClearFormats()
$new.Range("D2:G6").ClearFormats()

# --- 4) Prepend a "2022-Q1" summary row to the "总计" sheet -------------------
# NOTE: the sheet reference captured in $total before the Add() above tracks
# worksheet *position*, not identity - after inserting "2022-Q1" ahead of it,
# that old handle now resolves to the newly-added sheet instead. Re-fetch
# "总计" by name now that the sheet order/renaming has settled.
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Rows.Insert() leaves stray formatting on the new row's B:D cells (and no
# formatting at all on A2) - clear B2:D2 and then re-apply the bold/bordered
# "序号" look of column A from the row just below (the old row 2, now row 3).
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 1.73

# Re-number the "序号" (index) column A for the rows that got pushed down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Restore the originally-active sheet/tab (adding + renaming a sheet makes
# it the active one as a side effect).
$wb.Worksheets.Item("2021-Q2").Activate()

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
